$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 36, pushing the old row 36 down to row 38.
$ws.Rows.Item(36).Insert()
$ws.Rows.Item(36).Insert()

# --- Row 33: update date / quality / volume / prices ---
$ws.Cells.Item(33, 4).Value = 44985       # D33 Fecha
$ws.Cells.Item(33, 12).Value = "Primera"  # L33 Calidad
$ws.Cells.Item(33, 13).Value = 50         # M33 Volumen
$ws.Cells.Item(33, 14).Value = 3000       # N33 Precio minimo
$ws.Cells.Item(33, 15).Value = 3000       # O33 Precio maximo
$ws.Cells.Item(33, 16).Value = 3000       # P33 Precio promedio ponderado
$ws.Cells.Item(33, 19).Value = 1500       # S33 Precio $/Kg

# --- Row 34: update date / quality / volume / prices / origin ---
$ws.Cells.Item(34, 4).Value = 44985                      # D34 Fecha
$ws.Cells.Item(34, 12).Value = "Segunda"                 # L34 Calidad
$ws.Cells.Item(34, 13).Value = 50                         # M34 Volumen
$ws.Cells.Item(34, 14).Value = 2500                       # N34 Precio minimo
$ws.Cells.Item(34, 15).Value = 2500                       # O34 Precio maximo
$ws.Cells.Item(34, 16).Value = 2500                       # P34 Precio promedio ponderado
$ws.Cells.Item(34, 18).Value = "Provincia de Diguillín"   # R34 Origen
$ws.Cells.Item(34, 19).Value = 1250                       # S34 Precio $/Kg

# --- Row 35: update date / quality ---
$ws.Cells.Item(35, 4).Value = 44960    # D35 Fecha
$ws.Cells.Item(35, 12).Value = "Segunda" # L35 Calidad

# --- Row 36 (newly inserted): fill with the full record that used to be row 34 ---
$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(36, 3).Value = "Ñuble"
$ws.Cells.Item(36, 4).Value = 44174
$ws.Cells.Item(36, 5).Value = 16
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100101
$ws.Cells.Item(36, 8).Value = "Berries"
$ws.Cells.Item(36, 9).Value = 100101001
$ws.Cells.Item(36, 10).Value = "Arándano (blue)"
$ws.Cells.Item(36, 11).Value = "Sin especificar"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 150
$ws.Cells.Item(36, 14).Value = 3700
$ws.Cells.Item(36, 15).Value = 3800
$ws.Cells.Item(36, 16).Value = 3747
$ws.Cells.Item(36, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(36, 18).Value = "Provincia de Linares"
$ws.Cells.Item(36, 19).Value = 1874
$ws.Cells.Item(36, 20).Value = 2

# --- Row 37 (newly inserted): fill with the full record that used to be row 35 ---
$ws.Cells.Item(37, 1).Value = 7
$ws.Cells.Item(37, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(37, 3).Value = "Ñuble"
$ws.Cells.Item(37, 4).Value = 44944
$ws.Cells.Item(37, 5).Value = 16
$ws.Cells.Item(37, 6).Value = "Fruta"
$ws.Cells.Item(37, 7).Value = 100101
$ws.Cells.Item(37, 8).Value = "Berries"
$ws.Cells.Item(37, 9).Value = 100101001
$ws.Cells.Item(37, 10).Value = "Arándano (blue)"
$ws.Cells.Item(37, 11).Value = "Sin especificar"
$ws.Cells.Item(37, 12).Value = "Primera"
$ws.Cells.Item(37, 13).Value = 60
$ws.Cells.Item(37, 14).Value = 2500
$ws.Cells.Item(37, 15).Value = 2500
$ws.Cells.Item(37, 16).Value = 2500
$ws.Cells.Item(37, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(37, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(37, 19).Value = 1250
$ws.Cells.Item(37, 20).Value = 2

# Row 38 already holds the original row-36 record (shifted down by the insert),
# and its values are unchanged by this edit, so nothing further to do there.
